$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5833.5537
$ws.Range("I64").Value = 3403.2917
$ws.Range("J64").Value = 7656.25
$ws.Range("K64").Value = 3403.2917
$ws.Range("L64").Value = 7656.25
$ws.Range("M64").Value = -3155.2917
$ws.Range("N64").Value = -8152.25

$ws.Range("H67").Value = 5833.5537
$ws.Range("I67").Value = 3403.2917
$ws.Range("J67").Value = 7656.25
$ws.Range("K67").Value = 3403.2917
$ws.Range("L67").Value = 7656.25
$ws.Range("M67").Value = -2545.2917
$ws.Range("N67").Value = -9372.25

$ws.Range("H86").Value = 4310.846
$ws.Range("I86").Value = 4666.6665
$ws.Range("J86").Value = 4005.8572
$ws.Range("K86").Value = 4666.6665
$ws.Range("L86").Value = 4005.8572
$ws.Range("M86").Value = -3543.6665
$ws.Range("N86").Value = -6251.8572

$ws.Range("H89").Value = 4310.846
$ws.Range("I89").Value = 4666.6665
$ws.Range("J89").Value = 4005.8572
$ws.Range("K89").Value = 23333.3325
$ws.Range("L89").Value = 20029.286
$ws.Range("M89").Value = -17717.3325
$ws.Range("N89").Value = -31261.286

$ws.Range("H92").Value = 5813.1333
$ws.Range("I92").Value = 2585.9
$ws.Range("K92").Value = 2585.9
$ws.Range("M92").Value = -1337.9

$ws.Range("H129").Value = 1562.9333
$ws.Range("I129").Value = 767.63635
$ws.Range("K129").Value = 2302.90905
$ws.Range("M129").Value = 2697.09095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2672.0217
$ws.Range("I32").Value = 2720.2888
$ws.Range("K32").Value = 2720.2888
$ws.Range("M32").Value = -2433.2888

$ws.Range("H45").Value = 1112.5
$ws.Range("J45").Value = 1112.5
$ws.Range("L45").Value = 1112.5
$ws.Range("N45").Value = -1866.5

$ws.Range("H61").Value = 5443.4375
$ws.Range("I61").Value = 5443.4375
$ws.Range("K61").Value = 5443.4375
$ws.Range("M61").Value = -5231.4375

$ws.Range("H97").Value = 818.2
$ws.Range("I97").Value = 818.2
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 818.2
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -322.2
$ws.Range("N97").ClearContents()

$ws.Range("H110").Value = 1021.6667
$ws.Range("I110").Value = 1021.6667
$ws.Range("K110").Value = 1021.6667
$ws.Range("M110").Value = 1023.3333

$ws.Range("H132").Value = 4185.7827
$ws.Range("I132").Value = 2466.3845
$ws.Range("K132").Value = 7399.1535
$ws.Range("M132").Value = -4869.1535

$ws.Range("H136").Value = 5443.4375
$ws.Range("I136").Value = 5443.4375
$ws.Range("K136").Value = 16330.3125
$ws.Range("M136").Value = -13780.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 58789.5
$ws.Range("J2").Value = 58789.5
$ws.Range("L2").Value = 58789.5
$ws.Range("N2").Value = -59015.5

$ws.Range("H80").Value = 409.92
$ws.Range("I80").Value = 439.57144
$ws.Range("J80").Value = 398.3889
$ws.Range("K80").Value = 439.57144
$ws.Range("L80").Value = 398.3889
$ws.Range("M80").Value = 558.4285600000001
$ws.Range("N80").Value = -2394.3889

$ws.Range("H83").Value = 409.92
$ws.Range("I83").Value = 439.57144
$ws.Range("J83").Value = 398.3889
$ws.Range("K83").Value = 2197.8572
$ws.Range("L83").Value = 1991.9445
$ws.Range("M83").Value = 2794.1428
$ws.Range("N83").Value = -11975.9445

$ws.Range("H108").Value = 75000
$ws.Range("J108").Value = 75000
$ws.Range("L108").Value = 75000
$ws.Range("N108").Value = -82680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5122.4443
$ws.Range("I31").Value = 1992.8462
$ws.Range("K31").Value = 1992.8462
$ws.Range("M31").Value = -1697.8462

$ws.Range("H34").Value = 5122.4443
$ws.Range("I34").Value = 1992.8462
$ws.Range("K34").Value = 1992.8462
$ws.Range("M34").Value = -1790.8462

$ws.Range("H58").Value = 1345.6471
$ws.Range("I58").Value = 1521.909
$ws.Range("K58").Value = 1521.909
$ws.Range("M58").Value = -1318.909

$ws.Range("H86").Value = 1950
$ws.Range("I86").Value = 1900
$ws.Range("K86").Value = 1900
$ws.Range("M86").Value = -777

$ws.Range("H89").Value = 1950
$ws.Range("I89").Value = 1900
$ws.Range("K89").Value = 9500
$ws.Range("M89").Value = -3884

$ws.Range("H105").Value = 3055
$ws.Range("I105").Value = 3055
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3055
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1308
$ws.Range("N105").ClearContents()

$ws.Range("H136").Value = 1345.6471
$ws.Range("I136").Value = 1521.909
$ws.Range("K136").Value = 4565.727000000001
$ws.Range("M136").Value = -2015.727000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1410.8334
$ws.Range("I5").Value = 475.85715
$ws.Range("J5").Value = 2719.8
$ws.Range("K5").Value = 1427.57145
$ws.Range("L5").Value = 8159.400000000001
$ws.Range("M5").Value = -1315.57145
$ws.Range("N5").Value = -8383.400000000001

$ws.Range("H13").Value = 3144.4614
$ws.Range("I13").Value = 245
$ws.Range("K13").Value = 735
$ws.Range("M13").Value = -567

$ws.Range("H131").Value = 1881812.8
$ws.Range("J131").Value = 2901365.2
$ws.Range("L131").Value = 8704095.600000001
$ws.Range("N131").Value = -8714175.600000001

$ws.Range("H132").Value = 2512.4167
$ws.Range("I132").Value = 1562.25
$ws.Range("K132").Value = 14060.25
$ws.Range("M132").Value = -11530.25

$ws.Range("H135").Value = 1410.8334
$ws.Range("I135").Value = 475.85715
$ws.Range("J135").Value = 2719.8
$ws.Range("K135").Value = 4282.71435
$ws.Range("L135").Value = 24478.2
$ws.Range("M135").Value = -1747.71435
$ws.Range("N135").Value = -29548.2

$ws.Range("H140").Value = 64293.75
$ws.Range("I140").Value = 78370.766
$ws.Range("K140").Value = 235112.298
$ws.Range("M140").Value = -229932.298

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4551.9585
$ws.Range("I97").Value = 443.25
$ws.Range("J97").Value = 12769.375
$ws.Range("K97").Value = 443.25
$ws.Range("L97").Value = 12769.375
$ws.Range("M97").Value = 52.75
$ws.Range("N97").Value = -13761.375

$ws.Range("H102").Value = 1569.0476
$ws.Range("I102").Value = 1396.8
$ws.Range("K102").Value = 1396.8
$ws.Range("M102").Value = 225.2

$ws.Range("H132").Value = 2513.1177
$ws.Range("I132").Value = 2502.375
$ws.Range("K132").Value = 7507.125
$ws.Range("M132").Value = -4977.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 646.875
$ws.Range("I16").Value = 755
$ws.Range("K16").Value = 755
$ws.Range("M16").Value = -585

$ws.Range("H93").Value = 4993.5293
$ws.Range("I93").Value = 2148.8333
$ws.Range("J93").Value = 6545.1816
$ws.Range("K93").Value = 2148.8333
$ws.Range("L93").Value = 6545.1816
$ws.Range("M93").Value = -900.8332999999998
$ws.Range("N93").Value = -9041.1816

$ws.Range("H132").Value = 3386.889
$ws.Range("I132").Value = 3431.2856
$ws.Range("K132").Value = 10293.8568
$ws.Range("M132").Value = -7763.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws.Range("H122").Value = 3777.6365
$ws.Range("I122").Value = 1239.8572
$ws.Range("K122").Value = 3719.5716
$ws.Range("M122").Value = -1269.5716
